$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks (and their relationships) before rewriting the table
$ws.Cells.Hyperlinks.Delete()

# Clear out the old data range (keep header row 1 untouched)
$ws.Range("A2:C11").Clear()

# Write the refreshed, alphabetically sorted list of teams
$ws.Range("A2").Value = 'bugredasmissões'
$ws.Range("B2").Value = 19209079
$ws.Range("C2").Value = 'https://cartola.globo.com/#!/time/19209079'
$ws.Hyperlinks.Add($ws.Range("C2"), "https://cartola.globo.com/", '!/time/19209079') | Out-Null
$ws.Range("C2").Style = "Hyperlink"

$ws.Range("A3").Value = 'C R Juvenal'
$ws.Range("B3").Value = 1488983
$ws.Range("C3").Value = 'https://cartola.globo.com/#!/time/1488983'
$ws.Hyperlinks.Add($ws.Range("C3"), "https://cartola.globo.com/", '!/time/1488983') | Out-Null
$ws.Range("C3").Style = "Hyperlink"

$ws.Range("A4").Value = 'Doug Leal F.C'
$ws.Range("B4").Value = 287965
$ws.Range("C4").Value = 'https://cartola.globo.com/#!/time/287965'
$ws.Hyperlinks.Add($ws.Range("C4"), "https://cartola.globo.com/", '!/time/287965') | Out-Null
$ws.Range("C4").Style = "Hyperlink"

$ws.Range("A5").Value = 'Esquadrão Gazembrino'
$ws.Range("B5").Value = 2916559
$ws.Range("C5").Value = 'https://cartola.globo.com/#!/time/2916559'
$ws.Hyperlinks.Add($ws.Range("C5"), "https://cartola.globo.com/", '!/time/2916559') | Out-Null
$ws.Range("C5").Style = "Hyperlink"

$ws.Range("A6").Value = 'FBC Colorado'
$ws.Range("B6").Value = 186283
$ws.Range("C6").Value = 'https://cartola.globo.com/#!/time/186283'
$ws.Hyperlinks.Add($ws.Range("C6"), "https://cartola.globo.com/", '!/time/186283') | Out-Null
$ws.Range("C6").Style = "Hyperlink"

$ws.Range("A7").Value = 'GaúchoDaFronteira F.C'
$ws.Range("B7").Value = 2371918
$ws.Range("C7").Value = 'https://cartola.globo.com/#!/time/2371918'
$ws.Hyperlinks.Add($ws.Range("C7"), "https://cartola.globo.com/", '!/time/2371918') | Out-Null
$ws.Range("C7").Style = "Hyperlink"

$ws.Range("A8").Value = 'GE Bebum'
$ws.Range("B8").Value = 16411206
$ws.Range("C8").Value = 'https://cartola.globo.com/#!/time/16411206'
$ws.Hyperlinks.Add($ws.Range("C8"), "https://cartola.globo.com/", '!/time/16411206') | Out-Null
$ws.Range("C8").Style = "Hyperlink"

$ws.Range("A9").Value = 'Grêmio_Campeão_LA_27'
$ws.Range("B9").Value = 47775950
$ws.Range("C9").Value = 'https://cartola.globo.com/#!/time/47775950'
$ws.Hyperlinks.Add($ws.Range("C9"), "https://cartola.globo.com/", '!/time/47775950') | Out-Null
$ws.Range("C9").Style = "Hyperlink"

$ws.Range("A10").Value = 'JV5 Tricolor Gaúcho'
$ws.Range("B10").Value = 1747619
$ws.Range("C10").Value = 'https://cartola.globo.com/#!/time/1747619'
$ws.Hyperlinks.Add($ws.Range("C10"), "https://cartola.globo.com/", '!/time/1747619') | Out-Null
$ws.Range("C10").Style = "Hyperlink"

$ws.Range("A11").Value = 'La Primeira Patada Es Nuestra'
$ws.Range("B11").Value = 32966
$ws.Range("C11").Value = 'https://cartola.globo.com/#!/time/32966'
$ws.Hyperlinks.Add($ws.Range("C11"), "https://cartola.globo.com/", '!/time/32966') | Out-Null
$ws.Range("C11").Style = "Hyperlink"

$ws.Range("A12").Value = 'lsauer fc'
$ws.Range("B12").Value = 44810918
$ws.Range("C12").Value = 'https://cartola.globo.com/#!/time/44810918'
$ws.Hyperlinks.Add($ws.Range("C12"), "https://cartola.globo.com/", '!/time/44810918') | Out-Null
$ws.Range("C12").Style = "Hyperlink"

$ws.Range("A13").Value = 'Medonho´s F.C.'
$ws.Range("B13").Value = 1867254
$ws.Range("C13").Value = 'https://cartola.globo.com/#!/time/1867254'
$ws.Hyperlinks.Add($ws.Range("C13"), "https://cartola.globo.com/", '!/time/1867254') | Out-Null
$ws.Range("C13").Style = "Hyperlink"

$ws.Range("A14").Value = 'NHU PORÃ SAF.'
$ws.Range("B14").Value = 4088673
$ws.Range("C14").Value = 'https://cartola.globo.com/#!/time/4088673'
$ws.Hyperlinks.Add($ws.Range("C14"), "https://cartola.globo.com/", '!/time/4088673') | Out-Null
$ws.Range("C14").Style = "Hyperlink"

$ws.Range("A15").Value = 'Pontaç0 F.C.'
$ws.Range("B15").Value = 20651178
$ws.Range("C15").Value = 'https://cartola.globo.com/#!/time/20651178'
$ws.Hyperlinks.Add($ws.Range("C15"), "https://cartola.globo.com/", '!/time/20651178') | Out-Null
$ws.Range("C15").Style = "Hyperlink"

$ws.Range("A16").Value = 'SC 100 Sono'
$ws.Range("B16").Value = 14709358
$ws.Range("C16").Value = 'https://cartola.globo.com/#!/time/14709358'
$ws.Hyperlinks.Add($ws.Range("C16"), "https://cartola.globo.com/", '!/time/14709358') | Out-Null
$ws.Range("C16").Style = "Hyperlink"

$ws.Range("A17").Value = 'SC ÉoINTER!'
$ws.Range("B17").Value = 184499
$ws.Range("C17").Value = 'https://cartola.globo.com/#!/time/184499'
$ws.Hyperlinks.Add($ws.Range("C17"), "https://cartola.globo.com/", '!/time/184499') | Out-Null
$ws.Range("C17").Style = "Hyperlink"

$ws.Range("A18").Value = 'Texas Club 2026'
$ws.Range("B18").Value = 1273719
$ws.Range("C18").Value = 'https://cartola.globo.com/#!/time/1273719'
$ws.Hyperlinks.Add($ws.Range("C18"), "https://cartola.globo.com/", '!/time/1273719') | Out-Null
$ws.Range("C18").Style = "Hyperlink"
